# The workbook's "Sheet1" holds a token-tagging table (paragraph_id,
# sentence_number, word_index, word, tag). Rows 56-109 already duplicate the
# text of rows 2-55 (same sentences re-tagged), but column C ("word_index")
# in that block was stored as TEXT instead of a number. This edit:
#   1) fixes rows 56-109 so column C holds real numbers, and
#   2) appends a further duplicate of that block (rows 110-163), this time
#      keeping column C as TEXT (matching how rows 56-109 originally looked),
#      growing the sheet from A1:E109 to A1:E163.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcBlock  = $ws.Range("A56:E109")
$destBlock = $ws.Range("A110:E163")

# 1) Duplicate A56:E109 -> A110:E163 (values only, same layout/content).
$destBlock.Value = $srcBlock.Value2

# 2) The freshly duplicated column C (word_index) must stay TEXT, exactly
#    like rows 56-109 were before this edit. Force text storage, then strip
#    the temporary "@" formatting so no stray number-format style lingers.
$newC = $ws.Range("C110:C163")
$newC.NumberFormat = "@"
$newC.Value = $ws.Range("C56:C109").Value2
$newC.ClearFormats()

# 3) Now convert the original C56:C109 (word_index) from text to real numbers.
$ws.Range("C56:C109").Value = $ws.Range("C56:C109").Value2
